$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 7777
$ws.Range("A3").Value = 8888
$ws.Range("A4").Value = 9999
$ws.Range("A5").Value = 1111

$ws.Range("A2:C5").Select()
